$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows before row 104 to make room for the new common-setting API test cases,
# shifting the existing rows 104-142 down to 114-152.
$ws.Rows("104:113").Insert()

# Populate the newly inserted rows with the new test case data.
$ws.Cells.Item(104, 1).Value = "Admin_128"
$ws.Cells.Item(104, 2).Value = "@API Admin fetch the user-activities-filter"
$ws.Cells.Item(104, 3).Value = "passed"
$ws.Cells.Item(105, 1).Value = "Admin_129"
$ws.Cells.Item(105, 2).Value = "@API Endpoint validation for user-activities-filter"
$ws.Cells.Item(105, 3).Value = "passed"
$ws.Cells.Item(106, 1).Value = "Admin_131"
$ws.Cells.Item(106, 2).Value = "@API Admin_user-activities-list"
$ws.Cells.Item(106, 3).Value = "passed"
$ws.Cells.Item(107, 1).Value = "Admin_132"
$ws.Cells.Item(107, 2).Value = "@API Admin user-activities-list_validation of incorrect HTTP method"
$ws.Cells.Item(107, 3).Value = "passed"
$ws.Cells.Item(108, 1).Value = "Admin_133"
$ws.Cells.Item(108, 2).Value = "@API Admin user-activities-list_validation of invalid endpoint."
$ws.Cells.Item(108, 3).Value = "passed"
$ws.Cells.Item(109, 1).Value = "Admin_134"
$ws.Cells.Item(109, 2).Value = "@API Admin_user-activities-list-pagination"
$ws.Cells.Item(109, 3).Value = "passed"
$ws.Cells.Item(110, 1).Value = "Admin_135"
$ws.Cells.Item(110, 2).Value = "@API Admin user-activities-list-pagination_validation of incorrect HTTP method"
$ws.Cells.Item(110, 3).Value = "passed"
$ws.Cells.Item(111, 1).Value = "Admin_136"
$ws.Cells.Item(111, 2).Value = "@API Admin user-activities-list-pagination_validation of invalid endpoint."
$ws.Cells.Item(111, 3).Value = "passed"
$ws.Cells.Item(112, 1).Value = "Admin_0130"
$ws.Cells.Item(112, 2).Value = "@API Admin user-activities-filter-Header field validation - invalid"
$ws.Cells.Item(112, 3).Value = "failed"
$ws.Cells.Item(113, 1).Value = "AL_001T"
$ws.Cells.Item(113, 2).Value = "@API Admin Login Success with Mandatory Fields"
$ws.Cells.Item(113, 3).Value = "passed"

# Correct the status of the AL_001J login row (now shifted to row 133) from failed to passed.
$ws.Cells.Item(133, 3).Value = "passed"
